$d = $word.ActiveDocument

# 1. Ativação date change
$d.Content.Find.Execute(
    "Ativação: 01/01/2018",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ativação: 01/01/2025",
    2) | Out-Null

# 2. Objetivos (PT) - duplicate the sentence
$d.Content.Find.Execute(
    "Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas. Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas.",
    2) | Out-Null

# 3. Objetivos (EN, italic) - the empty italic run/paragraph gets new text
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Font.Italic -eq -1 -and $r.Text.Length -le 1) {
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = "Provide the student with practical knowledge in the technological processes of preparing fermented and distilled beverages. Provide the student with practical knowledge in the technological processes of preparing fermented and distilled beverages."
        break
    }
}

# 4. Programa resumido (PT)
$d.Content.Find.Execute(
    "Elaboração prática de cerveja, cachaça, fermentados e destilados de frutas, cereais e tuberculos, vinhos e análise sensorial.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Elaboração prática de cerveja, aguardente, licores e iogurtes.",
    2) | Out-Null

# 5. Programa resumido (EN, italic)
$d.Content.Find.Execute(
    "Practical elaboration of beer, cachaça, fermented and distilled of fruits, cereals and tubers, wines and sensorial analysis.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Practical classes for preparing beer, cachaça and yogurts.",
    2) | Out-Null

# 6. Programa (PT)
$d.Content.Find.Execute(
    "1. Elaboração de cerveja: matérias-primas, preparação do mosto, tecnologia de fermentação e maturação.2. Elaboração de aguardente: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.3. Elaboração de destilados de frutas: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.4. Elaboração e vinhos: matérias-primas, preparação do mosto, tecnologia de fermentação, maturação.5. Análise sensorial: teste sensorial das bebidas preparadas nos itens anteriores",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Elaboração de cerveja_ matérias-primas, preparação do mosto, tecnologia de fermentação e maturação. 2. Elaboração de aguardente_ matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação. 3. Elaboração de iogurte_ matérias-primas, preparação do leite, tecnologia de fermentação e acabamento_",
    2) | Out-Null

# 7. Programa (EN, italic)
$d.Content.Find.Execute(
    "1.Beer preparation; raw-materials; wort preparation; technology of fermentation and maturation.2.Spirits preparation: raw-materials; wort preparation; technology of fermentation, distillation, aging.3.Fruit distillates preparation: raw-materials; wort preparation; technology of fermentation, distillation, aging.4.Wines preparation: raw-materials, wort preparation, technology of fermentation, maturation.5.Sensorial analysis: sensorial test of the beverages prepared in the items above.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Beer production: raw materials, wort preparation, fermentation and maturation technology. 2. Production of cachaça: raw materials, must preparation, fermentation technology, distillation, maturation. 3. Yogurt production: raw materials, milk preparation, fermentation technology and finishing.",
    2) | Out-Null

# 8. Método (avaliação)
$d.Content.Find.Execute(
    "Relatórios e seminários sobre os experimentos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Avaliação teórica-prática, relatórios e/ou seminários sobre os experimentos.",
    2) | Out-Null

# 9. Critério (avaliação)
$d.Content.Find.Execute(
    "Média aritmética entre os relatórios e seminários",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Média aritimética entre a provas teórico-prática, relatórios e/ou seminários sobre os experimentos",
    2) | Out-Null

# 10. Bibliografia - old text contains a non-standard punctuation character that
# Find cannot reliably match, so replace the whole paragraph range instead.
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*AQUARONE*") {
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = "1) DA SILVA, N., JUNQUEIRA, V. C. A., DE ARRUDA SILVEIRA, N. F., TANIWAKI, M. H., GOMES, R. A. R., OKAZAKI, M. M. Manual de métodos de análise microbiológica de alimentos e água. Editora Blucher, 2017. 2) DA-SILVA, R.; LAGO-VANZELA, E. S.; BAFFI, M. A. Uvas e vinhos: química, bioquímica e microbiologia. São Paulo, Editora Senac, 2015. 3) DE OLIVEIRA MORAES, I. Biotecnologia Industrial: biotecnologia na produção de alimentos. Vol. 4. 2ª Ed. Editora Blucher, 2021. 4) MARTIN, J. G. P., DE DEA LINDNER, J. Microbiologia de alimentos fermentados. Editora Blucher, 2022. 5) MENEZES e SILVA, C.H.P. Microbiologia da cerveja - Do básico ao avançado, o guia definitivo. Editora LF, 2019. 6) MUXEL, A. A. Química da Cerveja: Uma Abordagem Química e Bioquímica das Matérias-Primas, Processo de Produção e da Composição dos Compostos de Sabores da Cerveja. Editora Appris, 2022. 7) VENTURINI FILHO, W. G. Bebidas alcoólicas: ciência e tecnologia. Vol. 1. Editora Blucher, 2021."
        break
    }
}

Write-Output "Done"
